$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.007.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.62%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.919.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.47%  '

# Row 4
$ws.Range("E4").Value = '  -0.26%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.64%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.43%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4597'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.45%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3825'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.59%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07732'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9812'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.94%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.25'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.36%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.933.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.92%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.963'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.16%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.686'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.25%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07026'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.58%  '

# Row 16
$ws.Range("E16").Value = '  -0.24%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.84%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009509'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.70%  '

# Row 20
$ws.Range("E20").Value = '  -0.44%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.034.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.81%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.336'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.62%  '

# Row 23
$ws.Range("E23").Value = '  +0.34%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.089'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.40%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.68%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.691'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.71%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.60%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.853'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.92%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09329'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '

# Row 31
$ws.Range("E31").Value = '  +1.33%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.114'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.80%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.255'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.047'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.80%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05705'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.07%  '

# Row 36
$ws.Range("E36").Value = '  +0.06%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.52%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02049'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.41%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.041'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.94%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.524'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.23%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5526'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.22%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1752'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.18%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000002983'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.381'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.35%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.214'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.38%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5191'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.33%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06901'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.96%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.782'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.34%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9996'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.50%  '
